$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AC2").Value = 11
$ws.Range("AE2").Value = 1000
$ws.Range("AH2").Value = 17
$ws.Range("AI2").Value = 80
$ws.Range("AK2").Value = 17
$ws.Range("AM2").Value = 70
$ws.Range("AN2").Value = 8.199999999999999
$ws.Range("AO2").Value = 80
$ws.Range("J2").Value = 4.1
$ws.Range("P2").Value = 2.48
$ws.Range("Q2").Value = 1.54
$ws.Range("S2").Value = 2.32
$ws.Range("Z2").Value = 40
# Row 3
$ws.Range("AA3").Value = 270
$ws.Range("AD3").Value = 30
$ws.Range("AE3").Value = 120
$ws.Range("AI3").Value = 100
$ws.Range("AJ3").Value = 12
$ws.Range("F3").Value = 1.44
$ws.Range("G3").Value = 1.45
$ws.Range("H3").Value = 8.800000000000001
$ws.Range("I3").Value = 9
$ws.Range("L3").Value = 1.32
$ws.Range("V3").Value = 1.12
# Row 4
$ws.Range("AA4").Value = 65
$ws.Range("AC4").Value = 9
$ws.Range("AD4").Value = 14.5
$ws.Range("AE4").Value = 36
$ws.Range("AI4").Value = 38
$ws.Range("AJ4").Value = 26
$ws.Range("AK4").Value = 19
$ws.Range("AN4").Value = 11
$ws.Range("AO4").Value = 26
$ws.Range("F4").Value = 2.1
$ws.Range("G4").Value = 2.12
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 3.65
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 4.1
$ws.Range("Q4").Value = 1.67
$ws.Range("R4").Value = 1.57
$ws.Range("T4").Value = 1.61
$ws.Range("U4").Value = 2.54
$ws.Range("V4").Value = 1.37
$ws.Range("W4").Value = 1.9
$ws.Range("Y4").Value = 18
$ws.Range("Z4").Value = 29
# Row 5
$ws.Range("AB5").Value = 16
$ws.Range("AD5").Value = 26
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 14.5
$ws.Range("AG5").Value = 12
$ws.Range("AH5").Value = 19.5
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 20
$ws.Range("AL5").Value = 26
$ws.Range("AN5").Value = 5.7
$ws.Range("AO5").Value = 55
$ws.Range("F5").Value = 1.48
$ws.Range("I5").Value = 6.8
$ws.Range("J5").Value = 4.4
$ws.Range("L5").Value = 1.21
$ws.Range("P5").Value = 2.62
$ws.Range("S5").Value = 1.89
$ws.Range("U5").Value = 2.36
$ws.Range("V5").Value = 1.17
$ws.Range("X5").Value = 38
# Row 6
$ws.Range("AA6").Value = 42
$ws.Range("AB6").Value = 20
$ws.Range("AC6").Value = 12
$ws.Range("AD6").Value = 15
$ws.Range("AE6").Value = 29
$ws.Range("AF6").Value = 28
$ws.Range("AG6").Value = 16
$ws.Range("AH6").Value = 18
$ws.Range("AI6").Value = 36
$ws.Range("AJ6").Value = 55
$ws.Range("AK6").Value = 34
$ws.Range("AL6").Value = 40
$ws.Range("AN6").Value = 21
$ws.Range("AO6").Value = 17
$ws.Range("F6").Value = 2.72
$ws.Range("J6").Value = 3.8
$ws.Range("P6").Value = 2.42
$ws.Range("R6").Value = 1.58
$ws.Range("U6").Value = 2.54
$ws.Range("V6").Value = 1.6
$ws.Range("X6").Value = 29
$ws.Range("Y6").Value = 18.5
$ws.Range("Z6").Value = 24
# Row 7
$ws.Range("R7").Value = 1.44
$ws.Range("S7").Value = 2.56
$ws.Range("T7").Value = 1.85
$ws.Range("V7").Value = 1.14
$ws.Range("Z7").Value = 1000
# Row 8
$ws.Range("G8").Value = 2.5
$ws.Range("I8").Value = 3.4
$ws.Range("M8").Value = 1.09
$ws.Range("O8").Value = 1.38
$ws.Range("Q8").Value = 2.18
$ws.Range("U8").Value = 2.06
$ws.Range("X8").Value = 11.5
# Row 9
$ws.Range("AA9").Value = 38
$ws.Range("AB9").Value = 19
$ws.Range("AC9").Value = 9.4
$ws.Range("AD9").Value = 12
$ws.Range("AG9").Value = 14
$ws.Range("AH9").Value = 14
$ws.Range("AI9").Value = 27
$ws.Range("AJ9").Value = 42
$ws.Range("AK9").Value = 25
$ws.Range("AM9").Value = 48
$ws.Range("AN9").Value = 14
$ws.Range("F9").Value = 2.82
$ws.Range("G9").Value = 2.84
$ws.Range("J9").Value = 3.9
$ws.Range("K9").Value = 3.95
$ws.Range("O9").Value = 1.18
$ws.Range("P9").Value = 2.76
$ws.Range("Q9").Value = 1.55
$ws.Range("S9").Value = 2.32
$ws.Range("U9").Value = 2.92
$ws.Range("X9").Value = 25
$ws.Range("Y9").Value = 17
# Row 10
$ws.Range("AJ10").Value = 340
$ws.Range("AL10").Value = 110
$ws.Range("F10").Value = 10
$ws.Range("G10").Value = 10.5
$ws.Range("Q10").Value = 1.63
$ws.Range("U10").Value = 2
$ws.Range("V10").Value = 3.6
$ws.Range("X10").Value = 24
# Row 11
$ws.Range("AD11").Value = 80
$ws.Range("AE11").Value = 500
$ws.Range("AI11").Value = 330
$ws.Range("Q11").Value = 1.49
$ws.Range("U11").Value = 1.71
$ws.Range("X11").Value = 34
$ws.Range("Y11").Value = 70
# Row 12
$ws.Range("AB12").Value = 15
$ws.Range("AI12").Value = 90
$ws.Range("AN12").Value = 3.45
$ws.Range("H12").Value = 11
$ws.Range("N12").Value = 8.6
$ws.Range("S12").Value = 1.93
$ws.Range("U12").Value = 2.28
$ws.Range("W12").Value = 4.2
# Row 13
$ws.Range("AA13").Value = 16
$ws.Range("AF13").Value = 46
$ws.Range("AH13").Value = 21
$ws.Range("AK13").Value = 80
$ws.Range("F13").Value = 5.9
$ws.Range("G13").Value = 6.2
$ws.Range("H13").Value = 1.66
$ws.Range("I13").Value = 1.68
$ws.Range("J13").Value = 4.3
$ws.Range("K13").Value = 4.4
$ws.Range("N13").Value = 4.3
$ws.Range("O13").Value = 1.28
$ws.Range("P13").Value = 2.12
$ws.Range("Q13").Value = 1.85
$ws.Range("R13").Value = 1.43
$ws.Range("T13").Value = 1.89
$ws.Range("U13").Value = 2.06
$ws.Range("V13").Value = 2.48
$ws.Range("X13").Value = 16.5
$ws.Range("Y13").Value = 9.199999999999999
# Row 14
$ws.Range("G14").Value = 3.15
$ws.Range("H14").Value = 2.42
$ws.Range("T14").Value = 1.62
$ws.Range("U14").Value = 2.56
# Row 15
$ws.Range("AB15").Value = 17.5
$ws.Range("AC15").Value = 11
$ws.Range("AE15").Value = 30
$ws.Range("AF15").Value = 23
$ws.Range("AG15").Value = 14
$ws.Range("AI15").Value = 36
$ws.Range("AJ15").Value = 42
$ws.Range("AK15").Value = 28
$ws.Range("AM15").Value = 1000
$ws.Range("AN15").Value = 16.5
$ws.Range("AO15").Value = 19
$ws.Range("H15").Value = 2.46
$ws.Range("I15").Value = 3.2
$ws.Range("K15").Value = 4.7
$ws.Range("L15").Value = 1.26
$ws.Range("M15").Value = 1.02
$ws.Range("N15").Value = 3.45
$ws.Range("O15").Value = 1.2
$ws.Range("P15").Value = 2.18
$ws.Range("R15").Value = 1.51
$ws.Range("T15").Value = 1.49
$ws.Range("U15").Value = 2.38
$ws.Range("V15").Value = 1.45
$ws.Range("Y15").Value = 18.5
# Row 16
$ws.Range("AB16").Value = 13
$ws.Range("F16").Value = 2.52
$ws.Range("G16").Value = 2.84
$ws.Range("H16").Value = 2.98
$ws.Range("M16").Value = 1.06
$ws.Range("W16").Value = 1.57
$ws.Range("Z16").Value = 25
